# Update countries & provincias Spain
#
# This updates the "Pais" sheet of paises.xlsx with refreshed COVID-19
# figures for several countries and bumps the "last updated" timestamp.
# A handful of rows are tied on "Casos totales" (column B) with their
# neighbours, so refreshing one country's numbers changes its sort rank
# and the country *labels* on adjacent rows effectively swap while the
# row that did not get new data keeps its old numbers one position down
# (or up). Those label swaps are encoded below together with the value
# changes for each affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = 'Datos actualizados a 19 de Septiembre de 2020 a las 09:19'

# --- Row 28: Ucrania (simple refresh, no reordering) -------------------
$ws.Range("B28").Value = 172712
$ws.Range("C28").Value = 3240
$ws.Range("D28").Value = 76754
$ws.Range("E28").Value = 92442
$ws.Range("G28").Value = 48
$ws.Range("H28").Value = 3516

# --- Row 62: Armenia (simple refresh, no reordering) --------------------
$ws.Range("B62").Value = 47154
$ws.Range("C62").Value = 244
$ws.Range("D62").Value = 42551
$ws.Range("E62").Value = 3675
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 928

# --- Rows 67/68: Azerbaiyan / Afganistan swap ranks ---------------------
$ws.Range("A67").Value = 'Afganistan'
$ws.Range("B67").Value = 38919
$ws.Range("C67").Value = 36
$ws.Range("D67").Value = 32576
$ws.Range("E67").Value = 4906
$ws.Range("H67").Value = 1437

$ws.Range("A68").Value = 'Azerbaiyan'
$ws.Range("B68").Value = 38894
$ws.Range("D68").Value = 36424
$ws.Range("E68").Value = 1898
$ws.Range("H68").Value = 572

# --- Row 77: Australia (simple refresh, no reordering) -------------------
$ws.Range("D77").Value = 23962
$ws.Range("E77").Value = 2079

# --- Rows 85/86: Republica de Macedonia / Hungria swap ranks -------------
$ws.Range("A85").Value = 'Hungria'
$ws.Range("B85").Value = 16920
$ws.Range("C85").Value = 809
$ws.Range("D85").Value = 4382
$ws.Range("E85").Value = 11863
$ws.Range("G85").Value = 6
$ws.Range("H85").Value = 675

$ws.Range("A86").Value = 'Republica de Macedonia'
$ws.Range("B86").Value = 16417
$ws.Range("D86").Value = 13732
$ws.Range("E86").Value = 2002
$ws.Range("H86").Value = 683

# --- Rows 139-142: Georgia moves up, Sri Lanka/Reunion/Bahamas shift down
$ws.Range("A139").Value = 'Georgia'
$ws.Range("B139").Value = 3306
$ws.Range("C139").Value = 187
$ws.Range("D139").Value = 1481
$ws.Range("E139").Value = 1806
$ws.Range("H139").Value = 19

$ws.Range("A140").Value = 'Sri Lanka'
$ws.Range("B140").Value = 3281
$ws.Range("D140").Value = 3060
$ws.Range("E140").Value = 208
$ws.Range("H140").Value = 13

$ws.Range("A141").Value = 'Reunion'
$ws.Range("B141").Value = 3194
$ws.Range("D141").Value = 1794
$ws.Range("E141").Value = 1385
$ws.Range("H141").Value = 15

$ws.Range("A142").Value = 'Bahamas'
$ws.Range("B142").Value = 3177
$ws.Range("D142").Value = 1626
$ws.Range("E142").Value = 1482
$ws.Range("H142").Value = 69

# --- Row 176: Taiwan (simple refresh, no reordering) ----------------------
$ws.Range("B176").Value = 506
$ws.Range("C176").Value = 3
$ws.Range("D176").Value = 479
$ws.Range("E176").Value = 20

# --- Rows 204/205: Santa Lucia / Timor Oriental swap ranks (tied values) -
$ws.Range("A204").Value = 'Timor Oriental'
$ws.Range("A205").Value = 'Santa Lucia'

# --- Rows 214/215: Montserrat / Islas Malvinas swap ranks -----------------
$ws.Range("A214").Value = 'Islas Malvinas'
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = 'Montserrat'
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
